# Update cryptos (24h snapshot) list per upstream GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.984.86'
$ws.Range("E2").Value = '  +1.20%  '

$ws.Range("D3").Value = '3.134.16'
$ws.Range("E3").Value = '  +1.81%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.39%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.29'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("E10").Value = '  +2.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.421'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.88%  '

$ws.Range("E12").Value = '  +2.92%  '

$ws.Range("D13").Value = '3.671.54'
$ws.Range("E13").Value = '  +1.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000168'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.45%  '

$ws.Range("D16").Value = '58.085.93'
$ws.Range("E16").Value = '  +1.25%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.93%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.133.03'
$ws.Range("E18").Value = '  +1.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.88%  '

$ws.Range("E25").Value = '  +3.32%  '

$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.65%  '

$ws.Range("D29").Value = '0.0₃0876'
$ws.Range("E29").Value = '  +1.38%  '

$ws.Range("E30").Value = '  +2.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.89%  '

$ws.Range("E33").Value = '  +6.25%  '

$ws.Range("E34").Value = '  +3.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.58%  '

$ws.Range("E37").Value = '  +9.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("E39").Value = '  +7.08%  '

$ws.Range("D40").Value = '2.641.92'
$ws.Range("E40").Value = '  +9.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0676'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.699'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.82%  '

$ws.Range("E45").Value = '  +4.45%  '

$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("E47").Value = '  +11.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.971'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.90%  '

$ws.Range("E50").Value = '  +3.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.749'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.84%  '
